# Apply scraped market-price updates to the Phoenix Profits workbook.
# Each worksheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) is a crafting-class Leve
# profit table; H:N are price/profit columns recomputed by the market-data
# runner. Cells whose HQ price info disappeared are cleared (not zeroed).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 13388.444
$ws.Range("I18").Value = 13388.444
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 13388.444
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -13104.444
$ws.Range("N18").ClearContents()

$ws.Range("H80").Value = 733.7895
$ws.Range("I80").Value = 593.5
$ws.Range("J80").Value = 798.53845
$ws.Range("K80").Value = 1780.5
$ws.Range("L80").Value = 2395.61535
$ws.Range("M80").Value = -782.5
$ws.Range("N80").Value = -4391.61535

$ws.Range("H83").Value = 733.7895
$ws.Range("I83").Value = 593.5
$ws.Range("J83").Value = 798.53845
$ws.Range("K83").Value = 5341.5
$ws.Range("L83").Value = 7186.84605
$ws.Range("M83").Value = -349.5
$ws.Range("N83").Value = -17170.84605

$ws.Range("H96").Value = 861.4
$ws.Range("I96").Value = 1019.2857
$ws.Range("J96").Value = 493
$ws.Range("K96").Value = 3057.8571
$ws.Range("L96").Value = 1479
$ws.Range("M96").Value = -1684.8571
$ws.Range("N96").Value = -4225

$ws.Range("H99").Value = 990
$ws.Range("I99").Value = 794.75
$ws.Range("K99").Value = 2384.25
$ws.Range("M99").Value = -886.25

$ws.Range("H106").Value = 4709.222
$ws.Range("I106").Value = 4990
$ws.Range("K106").Value = 4990
$ws.Range("M106").Value = -4359

$ws.Range("H107").Value = 6632
$ws.Range("I107").Value = 6144
$ws.Range("J107").Value = 8218
$ws.Range("K107").Value = 6144
$ws.Range("L107").Value = 8218
$ws.Range("M107").Value = -4224
$ws.Range("N107").Value = -12058

$ws.Range("H131").Value = 8368.643
$ws.Range("I131").Value = 3158
$ws.Range("J131").Value = 21395.25
$ws.Range("K131").Value = 9474
$ws.Range("L131").Value = 64185.75
$ws.Range("M131").Value = -4434
$ws.Range("N131").Value = -74265.75

$ws.Range("H132").Value = 3147.394
$ws.Range("I132").Value = 3124.6453
$ws.Range("K132").Value = 9373.9359
$ws.Range("M132").Value = -6843.9359

$ws.Range("H135").Value = 2290.4443
$ws.Range("I135").Value = 2352.5
$ws.Range("K135").Value = 21172.5
$ws.Range("M135").Value = -18637.5

$ws.Range("H138").Value = 3470.6
$ws.Range("I138").Value = 1846.3334
$ws.Range("J138").Value = 4445.16
$ws.Range("K138").Value = 5539.0002
$ws.Range("L138").Value = 13335.48
$ws.Range("M138").Value = -399.0002000000004
$ws.Range("N138").Value = -23615.48

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1683.4857
$ws.Range("I2").Value = 1488.1428
$ws.Range("K2").Value = 1488.1428
$ws.Range("M2").Value = -1375.1428

$ws.Range("H44").Value = 35199.8
$ws.Range("J44").Value = 35199.8
$ws.Range("L44").Value = 35199.8
$ws.Range("N44").Value = -36175.8

$ws.Range("H102").Value = 981.2414
$ws.Range("I102").Value = 917.9231
$ws.Range("K102").Value = 917.9231
$ws.Range("M102").Value = 704.0769

$ws.Range("H116").Value = 1683.4857
$ws.Range("I116").Value = 1488.1428
$ws.Range("K116").Value = 1488.1428
$ws.Range("M116").Value = 805.8571999999999

$ws.Range("H122").Value = 2632.0476
$ws.Range("I122").Value = 2593.3684
$ws.Range("K122").Value = 7780.1052
$ws.Range("M122").Value = -5330.1052

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1683.4857
$ws.Range("I3").Value = 1488.1428
$ws.Range("K3").Value = 1488.1428
$ws.Range("M3").Value = -1374.1428

$ws.Range("H94").Value = 819.7778
$ws.Range("I94").Value = 632.95
$ws.Range("J94").Value = 1353.5714
$ws.Range("K94").Value = 632.95
$ws.Range("L94").Value = 1353.5714
$ws.Range("M94").Value = -181.95
$ws.Range("N94").Value = -2255.5714

$ws.Range("H107").Value = 11715.929
$ws.Range("I107").Value = 13750.735
$ws.Range("K107").Value = 13750.735
$ws.Range("M107").Value = -11830.735

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 20498.75
$ws.Range("J26").Value = 20498.75
$ws.Range("L26").Value = 20498.75
$ws.Range("N26").Value = -21072.75

$ws.Range("H132").Value = 3723.4
$ws.Range("I132").Value = 3450.2307
$ws.Range("K132").Value = 10350.6921
$ws.Range("M132").Value = -7820.6921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 131.6923
$ws.Range("I6").Value = 131.6923
$ws.Range("K6").Value = 395.0769
$ws.Range("M6").Value = -282.0769

$ws.Range("H12").Value = 900.3871
$ws.Range("J12").Value = 1049.8077
$ws.Range("L12").Value = 3149.4231
$ws.Range("N12").Value = -3495.4231

$ws.Range("H14").Value = 2384.4546
$ws.Range("I14").Value = 2384.4546
$ws.Range("K14").Value = 7153.3638
$ws.Range("M14").Value = -6980.3638

$ws.Range("H29").Value = 11265.714
$ws.Range("I29").Value = 534.75
$ws.Range("K29").Value = 1604.25
$ws.Range("M29").Value = -1327.25

$ws.Range("H32").Value = 250
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H59").Value = 1201
$ws.Range("I59").Value = 601.3333
$ws.Range("K59").Value = 1803.9999
$ws.Range("M59").Value = -1263.9999

$ws.Range("H113").Value = 43480210
$ws.Range("J113").Value = 55557800
$ws.Range("L113").Value = 166673400
$ws.Range("N113").Value = -166677740

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3385.8447
$ws.Range("I132").Value = 3329.4468
$ws.Range("K132").Value = 9988.340400000001
$ws.Range("M132").Value = -7458.340400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 8345.68
$ws.Range("I61").Value = 8349.529
$ws.Range("K61").Value = 8349.529
$ws.Range("M61").Value = -8147.529

$ws.Range("H100").Value = 2451.5518
$ws.Range("J100").Value = 2450.4614
$ws.Range("L100").Value = 2450.4614
$ws.Range("N100").Value = -3532.4614

$ws.Range("H113").Value = 8345.68
$ws.Range("I113").Value = 8349.529
$ws.Range("K113").Value = 8349.529
$ws.Range("M113").Value = -6179.529

$ws.Range("H122").Value = 8129.0625
$ws.Range("I122").Value = 8321
$ws.Range("J122").Value = 5250
$ws.Range("K122").Value = 24963
$ws.Range("L122").Value = 15750
$ws.Range("M122").Value = -22513
$ws.Range("N122").Value = -20650

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 653.5897
$ws.Range("I113").Value = 659.27026
$ws.Range("J113").Value = 548.5
$ws.Range("K113").Value = 1977.81078
$ws.Range("L113").Value = 1645.5
$ws.Range("M113").Value = 192.18922
$ws.Range("N113").Value = -5985.5

$ws.Range("H122").Value = 7235.5576
$ws.Range("I122").Value = 7112
$ws.Range("K122").Value = 21336
$ws.Range("M122").Value = -18886
